$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.402.27"
$ws.Range("E2").Value = "  +3.07%  "

# Row 3
$ws.Range("D3").Value = "3.246.33"
$ws.Range("E3").Value = "  +6.06%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'577.14"
$ws.Range("E5").Value = "  +3.01%  "

# Row 6
$ws.Range("D6").Value = "'152.84"
$ws.Range("E6").Value = "  +7.13%  "

# Row 7
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.16%  "

# Row 8
$ws.Range("D8").Value = "3.238.24"
$ws.Range("E8").Value = "  +6.15%  "

# Row 9
$ws.Range("D9").Value = "'0.514"
$ws.Range("E9").Value = "  +4.05%  "

# Row 10
$ws.Range("D10").Value = "'7.09"
$ws.Range("E10").Value = "  +8.61%  "

# Row 11
$ws.Range("D11").Value = "'0.164"
$ws.Range("E11").Value = "  +4.14%  "

# Row 12
$ws.Range("D12").Value = "'0.488"
$ws.Range("E12").Value = "  +3.38%  "

# Row 13
$ws.Range("D13").Value = "'37.70"
$ws.Range("E13").Value = "  +2.32%  "

# Row 14
$ws.Range("D14").Value = "'0.0000234"
$ws.Range("E14").Value = "  +4.19%  "

# Row 15
$ws.Range("D15").Value = "3.765.90"
$ws.Range("E15").Value = "  +6.00%  "

# Row 16
$ws.Range("D16").Value = "'555.69"
$ws.Range("E16").Value = "  +11.31%  "

# Row 17
$ws.Range("D17").Value = "66.355.92"
$ws.Range("E17").Value = "  +2.91%  "

# Row 18
$ws.Range("D18").Value = "3.245.69"
$ws.Range("E18").Value = "  +5.77%  "

# Row 19
$ws.Range("E19").Value = "  +2.69%  "

# Row 20
$ws.Range("D20").Value = "'7.08"
$ws.Range("E20").Value = "  +4.60%  "

# Row 21
$ws.Range("D21").Value = "'14.37"
$ws.Range("E21").Value = "  +3.54%  "

# Row 22
$ws.Range("D22").Value = "'0.741"
$ws.Range("E22").Value = "  +6.47%  "

# Row 23
$ws.Range("D23").Value = "'7.76"
$ws.Range("E23").Value = "  +6.66%  "

# Row 24
$ws.Range("D24").Value = "'13.57"
$ws.Range("E24").Value = "  +4.90%  "

# Row 25
$ws.Range("D25").Value = "'81.74"
$ws.Range("E25").Value = "  +2.77%  "

# Row 26
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.48%  "

# Row 27
$ws.Range("D27").Value = "'9.27"
$ws.Range("E27").Value = "  +17.37%  "

# Row 28
$ws.Range("D28").Value = "'2.94"
$ws.Range("E28").Value = "  +5.49%  "

# Row 29
$ws.Range("D29").Value = "'2.23"
$ws.Range("E29").Value = "  +4.80%  "

# Row 30
$ws.Range("D30").Value = "'27.76"
$ws.Range("E30").Value = "  +5.31%  "

# Row 31
$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D31").Value = "'2.73"
$ws.Range("E31").Value = "  +2.11%  "

# Row 32
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.43%  "

# Row 33
$ws.Range("E33").Value = "  +4.70%  "

# Row 34
$ws.Range("D34").Value = "'560.11"
$ws.Range("E34").Value = "  +8.15%  "

# Row 35
$ws.Range("D35").Value = "'5.71"
$ws.Range("E35").Value = "  +2.44%  "

# Row 36
$ws.Range("D36").Value = "'6.37"
$ws.Range("E36").Value = "  +4.92%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.0457"
$ws.Range("E37").Value = "  +11.65%  "

# Row 38
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "'55.28"
$ws.Range("E38").Value = "  +3.60%  "

# Row 39
$ws.Range("D39").Value = "'0.0860"
$ws.Range("E39").Value = "  +5.97%  "

# Row 40
$ws.Range("D40").Value = "'0.131"
$ws.Range("E40").Value = "  +5.97%  "

# Row 41
$ws.Range("D41").Value = "'3.04"
$ws.Range("E41").Value = "  +12.90%  "

# Row 42
$ws.Range("D42").Value = "3.159.49"
$ws.Range("E42").Value = "  +7.17%  "

# Row 43
$ws.Range("D43").Value = "'8.58"
$ws.Range("E43").Value = "  +1.13%  "

# Row 44
$ws.Range("D44").Value = "'0.274"
$ws.Range("E44").Value = "  +9.86%  "

# Row 45
$ws.Range("D45").Value = "'2.30"
$ws.Range("E45").Value = "  +6.36%  "

# Row 46
$ws.Range("D46").Value = "'26.40"
$ws.Range("E46").Value = "  +3.29%  "

# Row 47
$ws.Range("E47").Value = "  +0.08%  "

# Row 48
$ws.Range("D48").Value = "0.0₃0553"
$ws.Range("E48").Value = "  +0.99%  "

# Row 49
$ws.Range("D49").Value = "'125.83"
$ws.Range("E49").Value = "  +3.97%  "

# Row 50
$ws.Range("D50").Value = "'0.113"
$ws.Range("E50").Value = "  +1.62%  "

# Row 51
$ws.Range("E51").Value = "  +6.21%  "
